$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the login rows (2-7) into their new positions -----------------
# (username / password / type-of-login), written row by row.
$ws.Range("A2").Value = "savitha.ip9@gmail.com"
$ws.Range("B2").Value = "Innovapath9"
$ws.Range("C2").Value = "FB"

$ws.Range("A3").Value = "savitha.ip9@gmail.com"
$ws.Range("B3").Value = "Innovapath1"
$ws.Range("C3").Value = "GP"

$ws.Range("A4").Value = "savitha.ip9@gmail.com"
$ws.Range("B4").Value = "Innovapath9"
$ws.Range("C4").Value = "LI"

$ws.Range("A5").Value = "SaviTalent"
$ws.Range("B5").Value = "Innovapath9"
$ws.Range("C5").Value = "GH"

$ws.Range("A6").Value = "raj@abc.com"
$ws.Range("B6").Value = "Saviraj1"
$ws.Range("C6").Value = "Normal"

$ws.Range("A7").Value = "savitha.ip9@gmail.com"
$ws.Range("B7").Value = "Innovapath9"
$ws.Range("C7").Value = "TW"

# --- Rebuild the hyperlinks so they follow their rows to the new spots ------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:savitha.ip9@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:savitha.ip9@gmail.com") | Out-Null

# A5's hyperlink carries a stale "display" (it shows raj@abc.com even though
# the cell itself still reads SaviTalent) - same mismatch the sheet already
# had on this row before it moved. Set TextToDisplay then restore the real
# cell text so the cached display text is kept without clobbering the value.
$hlA5 = $ws.Hyperlinks.Add($ws.Range("A5"), "mailto:raj@abc.com")
$hlA5.TextToDisplay = "raj@abc.com"
$ws.Range("A5").Value = "SaviTalent"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:savitha.ip9@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:savitha.ip9@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:raj@abc.com") | Out-Null

# Re-adding hyperlinks re-stamps the "Hyperlink" cell style; put the
# original style back so the cells keep referencing the same xf as before.
$ws.Range("A2:A7").Style = "Hyperlink"

# --- Update the sheet view: zoom + scroll position + new selection ---------
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 190
$ws.Range("A6:XFD6").Select() | Out-Null
